# day 43 - Java Summer
# Applies the UT11.docx edits described by the commit diff:
#  1. Remove the empty <w:lang/> run-property from the title run
#     " - Alertas y dialogos en Java" (best-effort; see note below).
#  2. Resize five tables (tblW / gridCol / tcW all need updating).
#  3. Split the "Int respuesta = ..." run after "Component " and drop a
#     _GoBack bookmark at the split point.
#  4. Remove the old _GoBack bookmark that used to sit further down,
#     just before "Después podemos tomar la respuesta...".
#  5. Drop the redundant <w:tblCellMar> override that lived inside the
#     5th table's <w:tblPrEx> (best-effort).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Table width updates.  dxa -> points is /20.  Each table needs its
#    PreferredWidth (tblW) AND its single column's Width (gridCol/tcW)
#    updated so all three XML numbers move together.
# ---------------------------------------------------------------------
$tableWidths = @(9819, 9943, 9943, 9958, 9912)
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    $w = $tableWidths[$i - 1]
    $pts = $w / 20.0
    $t.PreferredWidth = $pts
    $t.Columns.Item(1).Width = $pts
}

# ---------------------------------------------------------------------
# 2) Remove the bookmark that used to live right before "Después
#    podemos tomar la respuesta...".
# ---------------------------------------------------------------------
$old = $d.Bookmarks.Item("_GoBack")
if ($old -ne $null) {
    $old.Delete()
}

# ---------------------------------------------------------------------
# 3) Split the "Int respuesta ..." run and insert the _GoBack bookmark
#    right after "...showConfirmDialog(Component ".
# ---------------------------------------------------------------------
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("Int respuesta =JOptionPane. showConfirmDialog(Component ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitPoint = $d.Range($rng.End, $rng.End)
    $d.Bookmarks.Add("_GoBack", $splitPoint)
}

# ---------------------------------------------------------------------
# 4) Remove the " - Alertas y dialogos en Java" run's empty <w:lang/>.
# ---------------------------------------------------------------------
$rng2 = $d.Content.Duplicate
$rng2.Find.Execute(" - Alertas y dialogos en Java", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
